$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (rows 3-21) holds dates stored as literal text, e.g. "28/07/2022".
# They need to become "28-07-2022" (slashes -> dashes) while staying plain
# text. Some of them (day <= 12) are ambiguous and would otherwise be
# auto-parsed into date serials when written through .Value, so force the
# cell to Text format first, write the literal, then restore the cell's
# original (default/"Normal") style so no stray formatting is left behind.
$dates = @{
  3  = "28-07-2022"
  4  = "01-08-2022"
  5  = "04-08-2022"
  6  = "08-08-2022"
  7  = "11-08-2022"
  8  = "15-08-2022"
  9  = "18-08-2022"
  10 = "22-08-2022"
  11 = "25-08-2022"
  12 = "29-08-2022"
  13 = "01-09-2022"
  14 = "05-09-2022"
  15 = "08-09-2022"
  16 = "12-09-2022"
  17 = "15-09-2022"
  18 = "19-09-2022"
  19 = "22-09-2022"
  20 = "26-09-2022"
  21 = "29-09-2022"
}

foreach ($r in $dates.Keys) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$r]
    $cell.Style = "Normal"
}

# Row 3 attendance counts: Total Attendance Count (D) and Invalid (G) go
# from 0 to 1.
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

# Row 13 attendance counts: Total Attendance Count (D) and Real (E) go from
# 0 to 1, while Absent (H) goes from 1 to 0.
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("H13").Value = 0
